$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log row appended at row 15 (dimension grows from A1:F14 to A1:F15)
$ws.Range("A15").Value = "edit1"
$ws.Range("B15").Value = "riya-morankar"
$ws.Range("C15").Value = "Merged"
$ws.Range("D15").Value = "N/A"

# Force the Date column to stay plain text "2025-06-18" instead of being
# auto-converted into a date serial number by the Value setter.
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2025-06-18"
$ws.Range("E15").ClearFormats()

$ws.Range("F15").Value = "0298eb9d9af211aa8d48a882c377fe378d42a3d1"
